$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "74.724.42"
$ws.Range("E2").Value = "  +0.99%  "

$ws.Range("D3").Value = "2.811.23"
$ws.Range("E3").Value = "  +7.12%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "186.93"
$ws.Range("E5").Value = "  +0.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "590.82"
$ws.Range("E6").Value = "  +1.58%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.546"
$ws.Range("E8").Value = "  +2.87%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.188"
$ws.Range("E9").Value = "  -5.32%  "

$ws.Range("D10").Value = "2.810.28"
$ws.Range("E10").Value = "  +7.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.373"
$ws.Range("E11").Value = "  +4.37%  "

$ws.Range("E12").Value = "  -1.92%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.87"
$ws.Range("E13").Value = "  +4.24%  "

$ws.Range("D14").Value = "3.332.52"
$ws.Range("E14").Value = "  +7.39%  "

$ws.Range("D15").Value = "74.718.83"
$ws.Range("E15").Value = "  +1.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000185"
$ws.Range("E16").Value = "  -2.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.77"
$ws.Range("E17").Value = "  +1.50%  "

$ws.Range("D18").Value = "2.813.84"
$ws.Range("E18").Value = "  +7.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.03"
$ws.Range("E19").Value = "  -0.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.24"
$ws.Range("E20").Value = "  +3.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.06"
$ws.Range("E21").Value = "  +2.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.25"
$ws.Range("E22").Value = "  -2.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.07"
$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("E24").Value = "  +0.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.87"
$ws.Range("E25").Value = "  +1.53%  "

$ws.Range("D26").Value = "2.948.29"
$ws.Range("E26").Value = "  +6.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.13"
$ws.Range("E27").Value = "  +0.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.62"
$ws.Range("E28").Value = "  +3.07%  "

$ws.Range("E29").Value = "  +8.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "511.67"
$ws.Range("E31").Value = "  -1.91%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.39"
$ws.Range("E32").Value = "  +0.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.62"
$ws.Range("E33").Value = "  -0.42%  "

$ws.Range("E34").Value = "  +1.98%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.05"
$ws.Range("E36").Value = "  +0.94%  "

$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.88"
$ws.Range("E37").Value = "  +3.94%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.118"
$ws.Range("E38").Value = "  +0.07%  "

$ws.Range("E39").Value = "  +0.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "183.96"
$ws.Range("E40").Value = "  +13.79%  "

$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.339"
$ws.Range("E42").Value = "  +4.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.97"
$ws.Range("E43").Value = "  +1.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.65"
$ws.Range("E44").Value = "  -0.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.20"
$ws.Range("E45").Value = "  +2.39%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.01"
$ws.Range("E46").Value = "  +2.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0871"
$ws.Range("E47").Value = "  +2.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.31"
$ws.Range("E48").Value = "  -3.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.572"
$ws.Range("E49").Value = "  +9.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.70"
$ws.Range("E50").Value = "  +2.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.635"
$ws.Range("E51").Value = "  +8.50%  "
